$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.993.89"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.820.08"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.27"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4515"
$ws.Range("E7").Value = "  +6.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3699"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07284"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8559"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.72"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "1.800.85"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.644"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07105"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.325"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.11"
$ws.Range("E16").Value = "  +4.51%  "
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008787"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.97"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "26.897.30"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.56"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.226"
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.240"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.46"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08875"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.182"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7510"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.967"
$ws.Range("E33").Value = "  +5.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.435"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.099"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01965"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05242"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5295"
$ws.Range("E39").Value = "  +5.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.169"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.877"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1705"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5222"
$ws.Range("E43").Value = "  +10.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.519"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.66"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.978"
$ws.Range("E46").Value = "  +9.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.65"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.666"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9181"
$ws.Range("E51").Value = "  +0.40%  "
